$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.936.84"
$ws.Range("E2").Value = "  -1.36%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.065.20"
$ws.Range("E3").Value = "  -1.12%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "559.11"
$ws.Range("E5").Value = "  -0.42%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "142.96"
$ws.Range("E6").Value = "  -0.80%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.063.95"
$ws.Range("E8").Value = "  -1.03%  "
$ws.Range("E9").Value = "  +3.28%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.154"
$ws.Range("E10").Value = "  +0.90%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.20"
$ws.Range("E11").Value = "  -1.66%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.481"
$ws.Range("E12").Value = "  +1.56%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000233"
$ws.Range("E13").Value = "  +1.63%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.32"
$ws.Range("E14").Value = "  -0.28%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.566.96"
$ws.Range("E15").Value = "  -1.00%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.953.86"
$ws.Range("E16").Value = "  -1.36%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.063.36"
$ws.Range("E17").Value = "  -0.91%  "
$ws.Range("E18").Value = "  +0.13%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.79"
$ws.Range("E19").Value = "  -0.12%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "488.96"
$ws.Range("E20").Value = "  +1.75%  "
$ws.Range("E21").Value = "  +3.41%  "
$ws.Range("B22").Value = "Polygon"
$ws.Range("C22").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.687"
$ws.Range("E22").Value = "  -0.24%  "
$ws.Range("B23").Value = "InternetComputer(DFINITY)"
$ws.Range("C23").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.68"
$ws.Range("E23").Value = "  +8.48%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.55"
$ws.Range("E24").Value = "  -0.70%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "82.87"
$ws.Range("E25").Value = "  +1.90%  "
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("E27").Value = "  +0.72%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.15"
$ws.Range("E28").Value = "  -1.12%  "
$ws.Range("E29").Value = "  -0.78%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("E30").Value = "  +0.30%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "26.54"
$ws.Range("E31").Value = "  +1.29%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.16"
$ws.Range("E32").Value = "  +0.37%  "
$ws.Range("E33").Value = "  +1.83%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.71"
$ws.Range("E34").Value = "  +1.44%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.24"
$ws.Range("E35").Value = "  +1.16%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "55.15"
$ws.Range("E36").Value = "  +0.44%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0411"
$ws.Range("E37").Value = "  -0.18%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "445.17"
$ws.Range("E38").Value = "  -5.27%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0817"
$ws.Range("E39").Value = "  -2.54%  "
$ws.Range("B40").Value = "dogwifhat"
$ws.Range("C40").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.81"
$ws.Range("E40").Value = "  -6.33%  "
$ws.Range("B41").Value = "Maker"
$ws.Range("C41").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.030.13"
$ws.Range("E41").Value = "  +1.74%  "
$ws.Range("E42").Value = "  +0.87%  "
$ws.Range("E43").Value = "  +1.62%  "
$ws.Range("E44").Value = "  +5.21%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "27.87"
$ws.Range("E45").Value = "  -1.40%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.27"
$ws.Range("E46").Value = "  +4.99%  "
$ws.Range("E48").Value = "  +0.93%  "
$ws.Range("B49").Value = "PEPE"
$ws.Range("C49").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0₃0518"
$ws.Range("E49").Value = "  -1.61%  "
$ws.Range("B50").Value = "Monero"
$ws.Range("C50").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "118.14"
$ws.Range("E50").Value = "  +1.08%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.13"
$ws.Range("E51").Value = "  +2.13%  "
